# Updated symbol list on Tue Dec 27 02:57:50 UTC 2022 with GitHub Actions
#
# The "Price" column (D) is stored as text (numeric-looking strings such as
# "0.03300" that must keep trailing zeros), so every write to that column
# forces the cell to Text format first, writes the literal string, then
# restores the cell's style to "Normal" so no stray formatting is left
# behind (the cell keeps the same visual/style as before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

function Set-PlainValue {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $ws.Cells.Item($Row, $Col).Value = $Text
}

# --- Column D (Price) updates -------------------------------------------------
Set-TextValue 2  4 "243.71"
Set-TextValue 3  4 "23.24"
Set-TextValue 4  4 "5.412"
Set-TextValue 5  4 "0.05977"
Set-TextValue 6  4 "3.431"
Set-TextValue 7  4 "6.535"
Set-TextValue 8  4 "0.8112"
Set-TextValue 9  4 "0.9337"
Set-TextValue 11 4 "0.07421"
Set-TextValue 12 4 "0.03300"
Set-TextValue 13 4 "0.03070"
Set-TextValue 14 4 "0.09361"
Set-TextValue 15 4 "3.853"
Set-TextValue 16 4 "0.001579"
Set-TextValue 18 4 "0.0005932"
Set-TextValue 19 4 "0.005973"
Set-TextValue 20 4 "0.001279"
Set-TextValue 21 4 "0.004907"
Set-TextValue 22 4 "0.00006803"
Set-TextValue 23 4 "3.577"
Set-TextValue 40 4 "0.03969"
Set-TextValue 44 4 "0.009194"
Set-TextValue 45 4 "0.00005216"
Set-TextValue 47 4 "0.7252"
Set-TextValue 48 4 "0.002408"
Set-TextValue 49 4 "0.00002101"
Set-TextValue 50 4 "0.0002001"

# --- Row 20: Volume(1h) label gained a "Bestin24h" suffix ---------------------
Set-PlainValue 20 5 "19BitKanKANBestin24h"

# --- Rows 41-43: coin list reshuffled (CEJI / KickToken / BKEXToken) ---------
# Row 41 becomes BKEXToken
Set-PlainValue 41 2 "BKEXToken"
Set-PlainValue 41 3 "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue  41 4 "0.1081"
Set-PlainValue 41 5 "40BKEXTokenBKK"

# Row 42 becomes CEJI
Set-PlainValue 42 2 "CEJI"
Set-PlainValue 42 3 "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue  42 4 "0.002571"
Set-PlainValue 42 5 "41CEJICEJI"

# Row 43 becomes KickToken
Set-PlainValue 43 2 "KickToken"
Set-PlainValue 43 3 "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue  43 4 "0.003068"
Set-PlainValue 43 5 "42KickTokenKICKWorstin24h"

# --- Row 47: "Worstin24h" suffix removed ---------------------------------------
Set-PlainValue 47 5 "46CoinbaseStockTokenCOIN"
